## "starting the files for creating full data"
## Sheet "temp page for part no": insert a new "match group" column
## (with a running =B+100 formula) ahead of the existing "matched
## partno" column, add a new "entered overall ?" column with a couple
## of "yes" entries already filled in, remove the stale reviewer
## comment on A6, and tidy up the column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("temp page for part no")

# Remove the old "CHHS: email" reviewer note that was left on A6.
$comment = $ws.Range("A6").Comment
if ($comment) {
    $comment.Delete()
}

# Insert a new column C ("match group"), pushing the old "matched
# partno" column (and its data) from C to D.
$ws.Columns("C:C").Insert()

# New header row: A=Person, B=partno, C=match group, D=matched
# partno (unchanged data, just shifted right), E=entered overall ?
$ws.Range("E1").Value = "entered overall ?"
$ws.Range("C1").Value = "match group"

# Fill the new "match group" column with a running total formula
# (partno + 100) for every data row.
for ($r = 2; $r -le 57; $r++) {
    $ws.Range("C$r").Formula = "=B$r+100"
}

# A couple of rows already got their overall entered.
$ws.Range("E10").Value = "yes"
$ws.Range("E55").Value = "yes"

# Column sizing to fit the new content.
$ws.Columns("A:A").ColumnWidth = 18
$ws.Columns("E:E").ColumnWidth = 12.8

# Update the active selection to match where work left off.
$null = $ws.Range("D10").Select()
